# New PO forecast model
# Updates three sheets:
#  - "Weekly Quantity": append a new week row (41)
#  - "Monthly Trend": append a new month row (20)
#  - "PO Forecast": refresh forecast values for existing rows and extend the
#    series by one more forecast point (row 49)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Weekly Quantity  (append row 41)
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Cells.Item(41, 1).Value = 45662.99999999999
$wsWeekly.Cells.Item(41, 2).Value = 20
$wsWeekly.Cells.Item(41, 1).NumberFormat = $wsWeekly.Cells.Item(40, 1).NumberFormat

# ---------------------------------------------------------------------------
# Sheet: Monthly Trend  (append row 20)
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(20, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(20, 2).Value = 20
$wsMonthly.Cells.Item(20, 1).NumberFormat = $wsMonthly.Cells.Item(19, 1).NumberFormat

# ---------------------------------------------------------------------------
# Sheet: PO Forecast  (new forecast model -> rewrite dates/values for rows
# 2-49; rows 2-40 keep their original date, only the forecast qty changes;
# rows 41-49 get both a new date and a new forecast qty, with row 49 being a
# brand new row extending the series)
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$forecastRows = @(
    @(2, 45109.99999999999, 608),
    @(3, 45130.99999999999, 596),
    @(4, 45137.99999999999, 593),
    @(5, 45151.99999999999, 585),
    @(6, 45158.99999999999, 581),
    @(7, 45172.99999999999, 573),
    @(8, 45186.99999999999, 565),
    @(9, 45200.99999999999, 557),
    @(10, 45207.99999999999, 554),
    @(11, 45214.99999999999, 550),
    @(12, 45228.99999999999, 542),
    @(13, 45235.99999999999, 538),
    @(14, 45256.99999999999, 526),
    @(15, 45277.99999999999, 515),
    @(16, 45298.99999999999, 503),
    @(17, 45312.99999999999, 495),
    @(18, 45326.99999999999, 487),
    @(19, 45333.99999999999, 484),
    @(20, 45361.99999999999, 468),
    @(21, 45396.99999999999, 449),
    @(22, 45410.99999999999, 441),
    @(23, 45424.99999999999, 433),
    @(24, 45466.99999999999, 410),
    @(25, 45473.99999999999, 406),
    @(26, 45494.99999999999, 394),
    @(27, 45501.99999999999, 390),
    @(28, 45508.99999999999, 386),
    @(29, 45515.99999999999, 382),
    @(30, 45522.99999999999, 378),
    @(31, 45529.99999999999, 375),
    @(32, 45536.99999999999, 371),
    @(33, 45543.99999999999, 367),
    @(34, 45550.99999999999, 363),
    @(35, 45564.99999999999, 355),
    @(36, 45571.99999999999, 351),
    @(37, 45578.99999999999, 347),
    @(38, 45592.99999999999, 340),
    @(39, 45599.99999999999, 336),
    @(40, 45627.99999999999, 320),
    @(41, 45662.99999999999, 301),
    @(42, 45669.99999999999, 297),
    @(43, 45676.99999999999, 293),
    @(44, 45683.99999999999, 289),
    @(45, 45690.99999999999, 285),
    @(46, 45697.99999999999, 281),
    @(47, 45704.99999999999, 277),
    @(48, 45711.99999999999, 273),
    @(49, 45718.99999999999, 270)
)

$dateStyleSource = $wsForecast.Cells.Item(40, 1)

foreach ($row in $forecastRows) {
    $r = $row[0]
    $dateVal = $row[1]
    $qtyVal = $row[2]

    $wsForecast.Cells.Item($r, 1).Value = $dateVal
    $wsForecast.Cells.Item($r, 2).Value = $qtyVal
    $wsForecast.Cells.Item($r, 1).NumberFormat = $dateStyleSource.NumberFormat
}

Write-Output "PO forecast model updated"
